$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 175 is brand new: copy the date number format onto D175 only (matches D-column styling used elsewhere) ---
$ws.Range("D175").NumberFormat = $ws.Range("D174").NumberFormat

# --- Row 175 constant (template) columns, matching the fixed pattern used throughout the sheet ---
$ws.Range("A175").Value = 5
$ws.Range("B175").Value = 'Macroferia Regional de Talca'
$ws.Range("C175").Value = 'Maule'
$ws.Range("E175").Value = 7
$ws.Range("F175").Value = 100112008
$ws.Range("G175").Value = 'Coliflor'
$ws.Range("H175").Value = 'Sin especificar'
$ws.Range("N175").Value = '$/unidad'
$ws.Range("Q175").Value = 1
$ws.Range("R175").Value = 'Hortaliza'

# --- Row 48: brand new weekly record inserted at top of this block (everything below shifts down 1) ---
$ws.Range("D48").Value = 44525
$ws.Range("I48").Value = 'Primera'
$ws.Range("J48").Value = 5000
$ws.Range("K48").Value = 600
$ws.Range("L48").Value = 600
$ws.Range("M48").Value = 600
$ws.Range("O48").Value = 'Región del Maule'
$ws.Range("P48").Value = 600

# --- Rows 49-175: shift down by one (new row r = old row r-1) for D,I,J,K,L,M,O,P ---
$ws.Range("D49").Value = 44266
$ws.Range("I49").Value = 'Primera'
$ws.Range("J49").Value = 2200
$ws.Range("K49").Value = 700
$ws.Range("L49").Value = 800
$ws.Range("M49").Value = 709
$ws.Range("O49").Value = 'Región del Maule'
$ws.Range("P49").Value = 709
$ws.Range("D50").Value = 44433
$ws.Range("I50").Value = 'Primera'
$ws.Range("J50").Value = 2000
$ws.Range("K50").Value = 650
$ws.Range("L50").Value = 650
$ws.Range("M50").Value = 650
$ws.Range("O50").Value = 'Región del Maule'
$ws.Range("P50").Value = 650
$ws.Range("D51").Value = 44433
$ws.Range("I51").Value = 'Segunda'
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 500
$ws.Range("L51").Value = 500
$ws.Range("M51").Value = 500
$ws.Range("O51").Value = 'Región del Maule'
$ws.Range("P51").Value = 500
$ws.Range("D52").Value = 44159
$ws.Range("I52").Value = 'Primera'
$ws.Range("J52").Value = 2000
$ws.Range("K52").Value = 600
$ws.Range("L52").Value = 600
$ws.Range("M52").Value = 600
$ws.Range("O52").Value = 'Región del Maule'
$ws.Range("P52").Value = 600
$ws.Range("D53").Value = 44512
$ws.Range("I53").Value = 'Primera'
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 600
$ws.Range("L53").Value = 600
$ws.Range("M53").Value = 600
$ws.Range("O53").Value = 'Región del Maule'
$ws.Range("P53").Value = 600
$ws.Range("D54").Value = 44286
$ws.Range("I54").Value = 'Primera'
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 800
$ws.Range("L54").Value = 800
$ws.Range("M54").Value = 800
$ws.Range("O54").Value = 'Región del Maule'
$ws.Range("P54").Value = 800
$ws.Range("D55").Value = 44335
$ws.Range("I55").Value = 'Primera'
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 600
$ws.Range("M55").Value = 600
$ws.Range("O55").Value = 'Región del Maule'
$ws.Range("P55").Value = 600
$ws.Range("D56").Value = 44186
$ws.Range("I56").Value = 'Primera'
$ws.Range("J56").Value = 3000
$ws.Range("K56").Value = 500
$ws.Range("L56").Value = 500
$ws.Range("M56").Value = 500
$ws.Range("O56").Value = 'Región del Maule'
$ws.Range("P56").Value = 500
$ws.Range("D57").Value = 44460
$ws.Range("I57").Value = 'Primera'
$ws.Range("J57").Value = 3000
$ws.Range("K57").Value = 600
$ws.Range("L57").Value = 600
$ws.Range("M57").Value = 600
$ws.Range("O57").Value = 'Región del Maule'
$ws.Range("P57").Value = 600
$ws.Range("D58").Value = 44438
$ws.Range("I58").Value = 'Primera'
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 600
$ws.Range("L58").Value = 600
$ws.Range("M58").Value = 600
$ws.Range("O58").Value = 'Región del Maule'
$ws.Range("P58").Value = 600
$ws.Range("D59").Value = 44519
$ws.Range("I59").Value = 'Primera'
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 600
$ws.Range("L59").Value = 600
$ws.Range("M59").Value = 600
$ws.Range("O59").Value = 'Región del Maule'
$ws.Range("P59").Value = 600
$ws.Range("D60").Value = 44392
$ws.Range("I60").Value = 'Primera'
$ws.Range("J60").Value = 3000
$ws.Range("K60").Value = 700
$ws.Range("L60").Value = 700
$ws.Range("M60").Value = 700
$ws.Range("O60").Value = 'Región del Maule'
$ws.Range("P60").Value = 700
$ws.Range("D61").Value = 44355
$ws.Range("I61").Value = 'Primera'
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 450
$ws.Range("L61").Value = 450
$ws.Range("M61").Value = 450
$ws.Range("O61").Value = 'Región del Maule'
$ws.Range("P61").Value = 450
$ws.Range("D62").Value = 44489
$ws.Range("I62").Value = 'Primera'
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 600
$ws.Range("M62").Value = 600
$ws.Range("O62").Value = 'Región del Maule'
$ws.Range("P62").Value = 600
$ws.Range("D63").Value = 44434
$ws.Range("I63").Value = 'Primera'
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 600
$ws.Range("L63").Value = 600
$ws.Range("M63").Value = 600
$ws.Range("O63").Value = 'Región del Maule'
$ws.Range("P63").Value = 600
$ws.Range("D64").Value = 44434
$ws.Range("I64").Value = 'Segunda'
$ws.Range("J64").Value = 2000
$ws.Range("K64").Value = 500
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = 500
$ws.Range("O64").Value = 'Región del Maule'
$ws.Range("P64").Value = 500
$ws.Range("D65").Value = 44497
$ws.Range("I65").Value = 'Primera'
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 600
$ws.Range("L65").Value = 600
$ws.Range("M65").Value = 600
$ws.Range("O65").Value = 'Región del Maule'
$ws.Range("P65").Value = 600
$ws.Range("D66").Value = 44449
$ws.Range("I66").Value = 'Primera'
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 600
$ws.Range("L66").Value = 600
$ws.Range("M66").Value = 600
$ws.Range("O66").Value = 'Región del Maule'
$ws.Range("P66").Value = 600
$ws.Range("D67").Value = 44358
$ws.Range("I67").Value = 'Primera'
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 500
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 500
$ws.Range("O67").Value = 'Región del Maule'
$ws.Range("P67").Value = 500
$ws.Range("D68").Value = 44399
$ws.Range("I68").Value = 'Segunda'
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 500
$ws.Range("L68").Value = 500
$ws.Range("M68").Value = 500
$ws.Range("O68").Value = 'Región del Maule'
$ws.Range("P68").Value = 500
$ws.Range("D69").Value = 44298
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 700
$ws.Range("L69").Value = 700
$ws.Range("M69").Value = 700
$ws.Range("O69").Value = 'Región del Maule'
$ws.Range("P69").Value = 700
$ws.Range("D70").Value = 44405
$ws.Range("I70").Value = 'Segunda'
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 500
$ws.Range("L70").Value = 500
$ws.Range("M70").Value = 500
$ws.Range("O70").Value = 'Región del Maule'
$ws.Range("P70").Value = 500
$ws.Range("D71").Value = 44273
$ws.Range("I71").Value = 'Primera'
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 800
$ws.Range("L71").Value = 800
$ws.Range("M71").Value = 800
$ws.Range("O71").Value = 'Región del Maule'
$ws.Range("P71").Value = 800
$ws.Range("D72").Value = 44435
$ws.Range("I72").Value = 'Primera'
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 600
$ws.Range("L72").Value = 600
$ws.Range("M72").Value = 600
$ws.Range("O72").Value = 'Región Metropolitana'
$ws.Range("P72").Value = 600
$ws.Range("D73").Value = 44435
$ws.Range("I73").Value = 'Primera'
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 600
$ws.Range("L73").Value = 650
$ws.Range("M73").Value = 612
$ws.Range("O73").Value = 'Región del Maule'
$ws.Range("P73").Value = 612
$ws.Range("D74").Value = 44435
$ws.Range("I74").Value = 'Segunda'
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 500
$ws.Range("L74").Value = 500
$ws.Range("M74").Value = 500
$ws.Range("O74").Value = 'Región del Maule'
$ws.Range("P74").Value = 500
$ws.Range("D75").Value = 44328
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 300
$ws.Range("K75").Value = 700
$ws.Range("L75").Value = 700
$ws.Range("M75").Value = 700
$ws.Range("O75").Value = 'Región del Maule'
$ws.Range("P75").Value = 700
$ws.Range("D76").Value = 44277
$ws.Range("I76").Value = 'Primera'
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 800
$ws.Range("L76").Value = 800
$ws.Range("M76").Value = 800
$ws.Range("O76").Value = 'Región del Maule'
$ws.Range("P76").Value = 800
$ws.Range("D77").Value = 44516
$ws.Range("I77").Value = 'Primera'
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 500
$ws.Range("L77").Value = 500
$ws.Range("M77").Value = 500
$ws.Range("O77").Value = 'Región del Maule'
$ws.Range("P77").Value = 500
$ws.Range("D78").Value = 44168
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 600
$ws.Range("L78").Value = 600
$ws.Range("M78").Value = 600
$ws.Range("O78").Value = 'Región del Maule'
$ws.Range("P78").Value = 600
$ws.Range("D79").Value = 44475
$ws.Range("I79").Value = 'Primera'
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 600
$ws.Range("L79").Value = 600
$ws.Range("M79").Value = 600
$ws.Range("O79").Value = 'Región del Maule'
$ws.Range("P79").Value = 600
$ws.Range("D80").Value = 44419
$ws.Range("I80").Value = 'Segunda'
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 500
$ws.Range("L80").Value = 500
$ws.Range("M80").Value = 500
$ws.Range("O80").Value = 'Región del Maule'
$ws.Range("P80").Value = 500
$ws.Range("D81").Value = 44162
$ws.Range("I81").Value = 'Primera'
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 500
$ws.Range("L81").Value = 500
$ws.Range("M81").Value = 500
$ws.Range("O81").Value = 'Región del Maule'
$ws.Range("P81").Value = 500
$ws.Range("D82").Value = 44357
$ws.Range("I82").Value = 'Primera'
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 500
$ws.Range("L82").Value = 500
$ws.Range("M82").Value = 500
$ws.Range("O82").Value = 'Región del Maule'
$ws.Range("P82").Value = 500
$ws.Range("D83").Value = 44333
$ws.Range("I83").Value = 'Primera'
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 600
$ws.Range("L83").Value = 600
$ws.Range("M83").Value = 600
$ws.Range("O83").Value = 'Región del Maule'
$ws.Range("P83").Value = 600
$ws.Range("D84").Value = 44320
$ws.Range("I84").Value = 'Primera'
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 700
$ws.Range("L84").Value = 700
$ws.Range("M84").Value = 700
$ws.Range("O84").Value = 'Región del Maule'
$ws.Range("P84").Value = 700
$ws.Range("D85").Value = 44467
$ws.Range("I85").Value = 'Primera'
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 600
$ws.Range("L85").Value = 600
$ws.Range("M85").Value = 600
$ws.Range("O85").Value = 'Región del Maule'
$ws.Range("P85").Value = 600
$ws.Range("D86").Value = 44264
$ws.Range("I86").Value = 'Primera'
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 800
$ws.Range("L86").Value = 800
$ws.Range("M86").Value = 800
$ws.Range("O86").Value = 'Región del Maule'
$ws.Range("P86").Value = 800
$ws.Range("D87").Value = 44214
$ws.Range("I87").Value = 'Primera'
$ws.Range("J87").Value = 3000
$ws.Range("K87").Value = 700
$ws.Range("L87").Value = 700
$ws.Range("M87").Value = 700
$ws.Range("O87").Value = 'Región del Maule'
$ws.Range("P87").Value = 700
$ws.Range("D88").Value = 44167
$ws.Range("I88").Value = 'Primera'
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 600
$ws.Range("L88").Value = 600
$ws.Range("M88").Value = 600
$ws.Range("O88").Value = 'Región del Maule'
$ws.Range("P88").Value = 600
$ws.Range("D89").Value = 44291
$ws.Range("I89").Value = 'Primera'
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 700
$ws.Range("L89").Value = 700
$ws.Range("M89").Value = 700
$ws.Range("O89").Value = 'Región del Maule'
$ws.Range("P89").Value = 700
$ws.Range("D90").Value = 44293
$ws.Range("I90").Value = 'Primera'
$ws.Range("J90").Value = 3000
$ws.Range("K90").Value = 600
$ws.Range("L90").Value = 600
$ws.Range("M90").Value = 600
$ws.Range("O90").Value = 'Región del Maule'
$ws.Range("P90").Value = 600
$ws.Range("D91").Value = 44496
$ws.Range("I91").Value = 'Primera'
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 500
$ws.Range("L91").Value = 500
$ws.Range("M91").Value = 500
$ws.Range("O91").Value = 'Región del Maule'
$ws.Range("P91").Value = 500
$ws.Range("D92").Value = 44326
$ws.Range("I92").Value = 'Primera'
$ws.Range("J92").Value = 3000
$ws.Range("K92").Value = 500
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 500
$ws.Range("O92").Value = 'Región del Maule'
$ws.Range("P92").Value = 500
$ws.Range("D93").Value = 44302
$ws.Range("I93").Value = 'Primera'
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 600
$ws.Range("L93").Value = 600
$ws.Range("M93").Value = 600
$ws.Range("O93").Value = 'Región del Maule'
$ws.Range("P93").Value = 600
$ws.Range("D94").Value = 44292
$ws.Range("I94").Value = 'Primera'
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 700
$ws.Range("L94").Value = 700
$ws.Range("M94").Value = 700
$ws.Range("O94").Value = 'Región del Maule'
$ws.Range("P94").Value = 700
$ws.Range("D95").Value = 44308
$ws.Range("I95").Value = 'Primera'
$ws.Range("J95").Value = 3000
$ws.Range("K95").Value = 600
$ws.Range("L95").Value = 600
$ws.Range("M95").Value = 600
$ws.Range("O95").Value = 'Región del Maule'
$ws.Range("P95").Value = 600
$ws.Range("D96").Value = 44498
$ws.Range("I96").Value = 'Primera'
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 550
$ws.Range("L96").Value = 550
$ws.Range("M96").Value = 550
$ws.Range("O96").Value = 'Región del Maule'
$ws.Range("P96").Value = 550
$ws.Range("D97").Value = 44420
$ws.Range("I97").Value = 'Segunda'
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 500
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = 500
$ws.Range("O97").Value = 'Región del Maule'
$ws.Range("P97").Value = 500
$ws.Range("D98").Value = 44396
$ws.Range("I98").Value = 'Primera'
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 750
$ws.Range("M98").Value = 750
$ws.Range("O98").Value = 'Región Metropolitana'
$ws.Range("P98").Value = 750
$ws.Range("D99").Value = 44321
$ws.Range("I99").Value = 'Primera'
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 600
$ws.Range("O99").Value = 'Región del Maule'
$ws.Range("P99").Value = 600
$ws.Range("D100").Value = 44349
$ws.Range("I100").Value = 'Primera'
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 500
$ws.Range("L100").Value = 500
$ws.Range("M100").Value = 500
$ws.Range("O100").Value = 'Región del Maule'
$ws.Range("P100").Value = 500
$ws.Range("D101").Value = 44477
$ws.Range("I101").Value = 'Primera'
$ws.Range("J101").Value = 3000
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 600
$ws.Range("M101").Value = 600
$ws.Range("O101").Value = 'Región del Maule'
$ws.Range("P101").Value = 600
$ws.Range("D102").Value = 44487
$ws.Range("I102").Value = 'Primera'
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 600
$ws.Range("L102").Value = 600
$ws.Range("M102").Value = 600
$ws.Range("O102").Value = 'Región del Maule'
$ws.Range("P102").Value = 600
$ws.Range("D103").Value = 44452
$ws.Range("I103").Value = 'Primera'
$ws.Range("J103").Value = 3000
$ws.Range("K103").Value = 600
$ws.Range("L103").Value = 600
$ws.Range("M103").Value = 600
$ws.Range("O103").Value = 'Región del Maule'
$ws.Range("P103").Value = 600
$ws.Range("D104").Value = 44505
$ws.Range("I104").Value = 'Primera'
$ws.Range("J104").Value = 6000
$ws.Range("K104").Value = 500
$ws.Range("L104").Value = 500
$ws.Range("M104").Value = 500
$ws.Range("O104").Value = 'Región del Maule'
$ws.Range("P104").Value = 500
$ws.Range("D105").Value = 44306
$ws.Range("I105").Value = 'Primera'
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 600
$ws.Range("L105").Value = 600
$ws.Range("M105").Value = 600
$ws.Range("O105").Value = 'Región del Maule'
$ws.Range("P105").Value = 600
$ws.Range("D106").Value = 44509
$ws.Range("I106").Value = 'Primera'
$ws.Range("J106").Value = 5000
$ws.Range("K106").Value = 600
$ws.Range("L106").Value = 600
$ws.Range("M106").Value = 600
$ws.Range("O106").Value = 'Región del Maule'
$ws.Range("P106").Value = 600
$ws.Range("D107").Value = 44189
$ws.Range("I107").Value = 'Primera'
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 500
$ws.Range("O107").Value = 'Región del Maule'
$ws.Range("P107").Value = 500
$ws.Range("D108").Value = 44278
$ws.Range("I108").Value = 'Primera'
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 800
$ws.Range("L108").Value = 800
$ws.Range("M108").Value = 800
$ws.Range("O108").Value = 'Región del Maule'
$ws.Range("P108").Value = 800
$ws.Range("D109").Value = 44265
$ws.Range("I109").Value = 'Primera'
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 800
$ws.Range("L109").Value = 800
$ws.Range("M109").Value = 800
$ws.Range("O109").Value = 'Región del Maule'
$ws.Range("P109").Value = 800
$ws.Range("D110").Value = 44494
$ws.Range("I110").Value = 'Primera'
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 600
$ws.Range("L110").Value = 600
$ws.Range("M110").Value = 600
$ws.Range("O110").Value = 'Región del Maule'
$ws.Range("P110").Value = 600
$ws.Range("D111").Value = 44300
$ws.Range("I111").Value = 'Primera'
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 700
$ws.Range("L111").Value = 700
$ws.Range("M111").Value = 700
$ws.Range("O111").Value = 'Región del Maule'
$ws.Range("P111").Value = 700
$ws.Range("D112").Value = 44356
$ws.Range("I112").Value = 'Primera'
$ws.Range("J112").Value = 5000
$ws.Range("K112").Value = 450
$ws.Range("L112").Value = 450
$ws.Range("M112").Value = 450
$ws.Range("O112").Value = 'Región del Maule'
$ws.Range("P112").Value = 450
$ws.Range("D113").Value = 44469
$ws.Range("I113").Value = 'Primera'
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 500
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = 500
$ws.Range("O113").Value = 'Región del Maule'
$ws.Range("P113").Value = 500
$ws.Range("D114").Value = 44453
$ws.Range("I114").Value = 'Primera'
$ws.Range("J114").Value = 3000
$ws.Range("K114").Value = 600
$ws.Range("L114").Value = 600
$ws.Range("M114").Value = 600
$ws.Range("O114").Value = 'Región del Maule'
$ws.Range("P114").Value = 600
$ws.Range("D115").Value = 44518
$ws.Range("I115").Value = 'Primera'
$ws.Range("J115").Value = 4000
$ws.Range("K115").Value = 500
$ws.Range("L115").Value = 500
$ws.Range("M115").Value = 500
$ws.Range("O115").Value = 'Región del Maule'
$ws.Range("P115").Value = 500
$ws.Range("D116").Value = 44446
$ws.Range("I116").Value = 'Primera'
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 600
$ws.Range("L116").Value = 600
$ws.Range("M116").Value = 600
$ws.Range("O116").Value = 'Región del Maule'
$ws.Range("P116").Value = 600
$ws.Range("D117").Value = 44463
$ws.Range("I117").Value = 'Primera'
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 600
$ws.Range("L117").Value = 600
$ws.Range("M117").Value = 600
$ws.Range("O117").Value = 'Región del Maule'
$ws.Range("P117").Value = 600
$ws.Range("D118").Value = 44323
$ws.Range("I118").Value = 'Primera'
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 600
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 600
$ws.Range("O118").Value = 'Región del Maule'
$ws.Range("P118").Value = 600
$ws.Range("D119").Value = 44417
$ws.Range("I119").Value = 'Segunda'
$ws.Range("J119").Value = 3000
$ws.Range("K119").Value = 500
$ws.Range("L119").Value = 500
$ws.Range("M119").Value = 500
$ws.Range("O119").Value = 'Región del Maule'
$ws.Range("P119").Value = 500
$ws.Range("D120").Value = 44445
$ws.Range("I120").Value = 'Primera'
$ws.Range("J120").Value = 4000
$ws.Range("K120").Value = 600
$ws.Range("L120").Value = 600
$ws.Range("M120").Value = 600
$ws.Range("O120").Value = 'Región del Maule'
$ws.Range("P120").Value = 600
$ws.Range("D121").Value = 44342
$ws.Range("I121").Value = 'Primera'
$ws.Range("J121").Value = 4000
$ws.Range("K121").Value = 700
$ws.Range("L121").Value = 700
$ws.Range("M121").Value = 700
$ws.Range("O121").Value = 'Región del Maule'
$ws.Range("P121").Value = 700
$ws.Range("D122").Value = 44523
$ws.Range("I122").Value = 'Primera'
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 600
$ws.Range("L122").Value = 600
$ws.Range("M122").Value = 600
$ws.Range("O122").Value = 'Región del Maule'
$ws.Range("P122").Value = 600
$ws.Range("D123").Value = 44406
$ws.Range("I123").Value = 'Segunda'
$ws.Range("J123").Value = 5000
$ws.Range("K123").Value = 400
$ws.Range("L123").Value = 400
$ws.Range("M123").Value = 400
$ws.Range("O123").Value = 'Región del Maule'
$ws.Range("P123").Value = 400
$ws.Range("D124").Value = 44295
$ws.Range("I124").Value = 'Primera'
$ws.Range("J124").Value = 4000
$ws.Range("K124").Value = 700
$ws.Range("L124").Value = 700
$ws.Range("M124").Value = 700
$ws.Range("O124").Value = 'Región del Maule'
$ws.Range("P124").Value = 700
$ws.Range("D125").Value = 44270
$ws.Range("I125").Value = 'Primera'
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 800
$ws.Range("L125").Value = 800
$ws.Range("M125").Value = 800
$ws.Range("O125").Value = 'Región del Maule'
$ws.Range("P125").Value = 800
$ws.Range("D126").Value = 44363
$ws.Range("I126").Value = 'Primera'
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 400
$ws.Range("L126").Value = 400
$ws.Range("M126").Value = 400
$ws.Range("O126").Value = 'Región del Maule'
$ws.Range("P126").Value = 400
$ws.Range("D127").Value = 44299
$ws.Range("I127").Value = 'Primera'
$ws.Range("J127").Value = 4000
$ws.Range("K127").Value = 700
$ws.Range("L127").Value = 700
$ws.Range("M127").Value = 700
$ws.Range("O127").Value = 'Región del Maule'
$ws.Range("P127").Value = 700
$ws.Range("D128").Value = 44336
$ws.Range("I128").Value = 'Primera'
$ws.Range("J128").Value = 3000
$ws.Range("K128").Value = 600
$ws.Range("L128").Value = 600
$ws.Range("M128").Value = 600
$ws.Range("O128").Value = 'Región del Maule'
$ws.Range("P128").Value = 600
$ws.Range("D129").Value = 44372
$ws.Range("I129").Value = 'Primera'
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 400
$ws.Range("L129").Value = 400
$ws.Range("M129").Value = 400
$ws.Range("O129").Value = 'Región del Maule'
$ws.Range("P129").Value = 400
$ws.Range("D130").Value = 44403
$ws.Range("I130").Value = 'Primera'
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 700
$ws.Range("L130").Value = 700
$ws.Range("M130").Value = 700
$ws.Range("O130").Value = 'Región Metropolitana'
$ws.Range("P130").Value = 700
$ws.Range("D131").Value = 44403
$ws.Range("I131").Value = 'Segunda'
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 500
$ws.Range("L131").Value = 500
$ws.Range("M131").Value = 500
$ws.Range("O131").Value = 'Región del Maule'
$ws.Range("P131").Value = 500
$ws.Range("D132").Value = 44169
$ws.Range("I132").Value = 'Primera'
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 600
$ws.Range("L132").Value = 600
$ws.Range("M132").Value = 600
$ws.Range("O132").Value = 'Región del Maule'
$ws.Range("P132").Value = 600
$ws.Range("D133").Value = 44376
$ws.Range("I133").Value = 'Primera'
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 400
$ws.Range("L133").Value = 400
$ws.Range("M133").Value = 400
$ws.Range("O133").Value = 'Región del Maule'
$ws.Range("P133").Value = 400
$ws.Range("D134").Value = 44172
$ws.Range("I134").Value = 'Primera'
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 600
$ws.Range("L134").Value = 600
$ws.Range("M134").Value = 600
$ws.Range("O134").Value = 'Región del Maule'
$ws.Range("P134").Value = 600
$ws.Range("D135").Value = 44421
$ws.Range("I135").Value = 'Segunda'
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 500
$ws.Range("L135").Value = 500
$ws.Range("M135").Value = 500
$ws.Range("O135").Value = 'Región del Maule'
$ws.Range("P135").Value = 500
$ws.Range("D136").Value = 44431
$ws.Range("I136").Value = 'Primera'
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 600
$ws.Range("L136").Value = 600
$ws.Range("M136").Value = 600
$ws.Range("O136").Value = 'Región Metropolitana'
$ws.Range("P136").Value = 600
$ws.Range("D137").Value = 44426
$ws.Range("I137").Value = 'Segunda'
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 500
$ws.Range("L137").Value = 500
$ws.Range("M137").Value = 500
$ws.Range("O137").Value = 'Región del Maule'
$ws.Range("P137").Value = 500
$ws.Range("D138").Value = 44448
$ws.Range("I138").Value = 'Primera'
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 600
$ws.Range("L138").Value = 600
$ws.Range("M138").Value = 600
$ws.Range("O138").Value = 'Región del Maule'
$ws.Range("P138").Value = 600
$ws.Range("D139").Value = 44362
$ws.Range("I139").Value = 'Primera'
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 400
$ws.Range("L139").Value = 400
$ws.Range("M139").Value = 400
$ws.Range("O139").Value = 'Región del Maule'
$ws.Range("P139").Value = 400
$ws.Range("D140").Value = 44176
$ws.Range("I140").Value = 'Primera'
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 600
$ws.Range("L140").Value = 600
$ws.Range("M140").Value = 600
$ws.Range("O140").Value = 'Región del Maule'
$ws.Range("P140").Value = 600
$ws.Range("D141").Value = 44301
$ws.Range("I141").Value = 'Primera'
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 600
$ws.Range("L141").Value = 600
$ws.Range("M141").Value = 600
$ws.Range("O141").Value = 'Región del Maule'
$ws.Range("P141").Value = 600
$ws.Range("D142").Value = 44407
$ws.Range("I142").Value = 'Segunda'
$ws.Range("J142").Value = 3000
$ws.Range("K142").Value = 500
$ws.Range("L142").Value = 500
$ws.Range("M142").Value = 500
$ws.Range("O142").Value = 'Región del Maule'
$ws.Range("P142").Value = 500
$ws.Range("D143").Value = 44284
$ws.Range("I143").Value = 'Primera'
$ws.Range("J143").Value = 5000
$ws.Range("K143").Value = 700
$ws.Range("L143").Value = 800
$ws.Range("M143").Value = 740
$ws.Range("O143").Value = 'Región del Maule'
$ws.Range("P143").Value = 740
$ws.Range("D144").Value = 44441
$ws.Range("I144").Value = 'Primera'
$ws.Range("J144").Value = 3000
$ws.Range("K144").Value = 600
$ws.Range("L144").Value = 600
$ws.Range("M144").Value = 600
$ws.Range("O144").Value = 'Región del Maule'
$ws.Range("P144").Value = 600
$ws.Range("D145").Value = 44279
$ws.Range("I145").Value = 'Primera'
$ws.Range("J145").Value = 3000
$ws.Range("K145").Value = 800
$ws.Range("L145").Value = 800
$ws.Range("M145").Value = 800
$ws.Range("O145").Value = 'Región del Maule'
$ws.Range("P145").Value = 800
$ws.Range("D146").Value = 44341
$ws.Range("I146").Value = 'Primera'
$ws.Range("J146").Value = 3000
$ws.Range("K146").Value = 700
$ws.Range("L146").Value = 700
$ws.Range("M146").Value = 700
$ws.Range("O146").Value = 'Región del Maule'
$ws.Range("P146").Value = 700
$ws.Range("D147").Value = 44504
$ws.Range("I147").Value = 'Primera'
$ws.Range("J147").Value = 5000
$ws.Range("K147").Value = 600
$ws.Range("L147").Value = 600
$ws.Range("M147").Value = 600
$ws.Range("O147").Value = 'Región del Maule'
$ws.Range("P147").Value = 600
$ws.Range("D148").Value = 44350
$ws.Range("I148").Value = 'Primera'
$ws.Range("J148").Value = 5000
$ws.Range("K148").Value = 500
$ws.Range("L148").Value = 500
$ws.Range("M148").Value = 500
$ws.Range("O148").Value = 'Región del Maule'
$ws.Range("P148").Value = 500
$ws.Range("D149").Value = 44312
$ws.Range("I149").Value = 'Primera'
$ws.Range("J149").Value = 5000
$ws.Range("K149").Value = 600
$ws.Range("L149").Value = 600
$ws.Range("M149").Value = 600
$ws.Range("O149").Value = 'Provincia de Chacabuco'
$ws.Range("P149").Value = 600
$ws.Range("D150").Value = 44384
$ws.Range("I150").Value = 'Primera'
$ws.Range("J150").Value = 4000
$ws.Range("K150").Value = 600
$ws.Range("L150").Value = 600
$ws.Range("M150").Value = 600
$ws.Range("O150").Value = 'Región del Maule'
$ws.Range("P150").Value = 600
$ws.Range("D151").Value = 44329
$ws.Range("I151").Value = 'Primera'
$ws.Range("J151").Value = 3000
$ws.Range("K151").Value = 650
$ws.Range("L151").Value = 650
$ws.Range("M151").Value = 650
$ws.Range("O151").Value = 'Región del Maule'
$ws.Range("P151").Value = 650
$ws.Range("D152").Value = 44522
$ws.Range("I152").Value = 'Primera'
$ws.Range("J152").Value = 4000
$ws.Range("K152").Value = 600
$ws.Range("L152").Value = 600
$ws.Range("M152").Value = 600
$ws.Range("O152").Value = 'Región del Maule'
$ws.Range("P152").Value = 600
$ws.Range("D153").Value = 44491
$ws.Range("I153").Value = 'Primera'
$ws.Range("J153").Value = 3000
$ws.Range("K153").Value = 700
$ws.Range("L153").Value = 700
$ws.Range("M153").Value = 700
$ws.Range("O153").Value = 'Región del Maule'
$ws.Range("P153").Value = 700
$ws.Range("D154").Value = 44272
$ws.Range("I154").Value = 'Primera'
$ws.Range("J154").Value = 3000
$ws.Range("K154").Value = 800
$ws.Range("L154").Value = 800
$ws.Range("M154").Value = 800
$ws.Range("O154").Value = 'Región del Maule'
$ws.Range("P154").Value = 800
$ws.Range("D155").Value = 44305
$ws.Range("I155").Value = 'Primera'
$ws.Range("J155").Value = 4000
$ws.Range("K155").Value = 600
$ws.Range("L155").Value = 600
$ws.Range("M155").Value = 600
$ws.Range("O155").Value = 'Región del Maule'
$ws.Range("P155").Value = 600
$ws.Range("D156").Value = 44166
$ws.Range("I156").Value = 'Primera'
$ws.Range("J156").Value = 3000
$ws.Range("K156").Value = 700
$ws.Range("L156").Value = 700
$ws.Range("M156").Value = 700
$ws.Range("O156").Value = 'Región del Maule'
$ws.Range("P156").Value = 700
$ws.Range("D157").Value = 44315
$ws.Range("I157").Value = 'Primera'
$ws.Range("J157").Value = 3000
$ws.Range("K157").Value = 600
$ws.Range("L157").Value = 600
$ws.Range("M157").Value = 600
$ws.Range("O157").Value = 'Región del Maule'
$ws.Range("P157").Value = 600
$ws.Range("D158").Value = 44348
$ws.Range("I158").Value = 'Primera'
$ws.Range("J158").Value = 6000
$ws.Range("K158").Value = 450
$ws.Range("L158").Value = 450
$ws.Range("M158").Value = 450
$ws.Range("O158").Value = 'Región del Maule'
$ws.Range("P158").Value = 450
$ws.Range("D159").Value = 44322
$ws.Range("I159").Value = 'Primera'
$ws.Range("J159").Value = 3000
$ws.Range("K159").Value = 600
$ws.Range("L159").Value = 600
$ws.Range("M159").Value = 600
$ws.Range("O159").Value = 'Región del Maule'
$ws.Range("P159").Value = 600
$ws.Range("D160").Value = 44495
$ws.Range("I160").Value = 'Primera'
$ws.Range("J160").Value = 5000
$ws.Range("K160").Value = 500
$ws.Range("L160").Value = 500
$ws.Range("M160").Value = 500
$ws.Range("O160").Value = 'Región del Maule'
$ws.Range("P160").Value = 500
$ws.Range("D161").Value = 44327
$ws.Range("I161").Value = 'Primera'
$ws.Range("J161").Value = 3000
$ws.Range("K161").Value = 700
$ws.Range("L161").Value = 700
$ws.Range("M161").Value = 700
$ws.Range("O161").Value = 'Región del Maule'
$ws.Range("P161").Value = 700
$ws.Range("D162").Value = 44510
$ws.Range("I162").Value = 'Primera'
$ws.Range("J162").Value = 4000
$ws.Range("K162").Value = 600
$ws.Range("L162").Value = 600
$ws.Range("M162").Value = 600
$ws.Range("O162").Value = 'Región del Maule'
$ws.Range("P162").Value = 600
$ws.Range("D163").Value = 44161
$ws.Range("I163").Value = 'Primera'
$ws.Range("J163").Value = 4000
$ws.Range("K163").Value = 500
$ws.Range("L163").Value = 500
$ws.Range("M163").Value = 500
$ws.Range("O163").Value = 'Región del Maule'
$ws.Range("P163").Value = 500
$ws.Range("D164").Value = 44517
$ws.Range("I164").Value = 'Primera'
$ws.Range("J164").Value = 5000
$ws.Range("K164").Value = 500
$ws.Range("L164").Value = 500
$ws.Range("M164").Value = 500
$ws.Range("O164").Value = 'Región del Maule'
$ws.Range("P164").Value = 500
$ws.Range("D165").Value = 44391
$ws.Range("I165").Value = 'Primera'
$ws.Range("J165").Value = 3000
$ws.Range("K165").Value = 700
$ws.Range("L165").Value = 700
$ws.Range("M165").Value = 700
$ws.Range("O165").Value = 'Región del Maule'
$ws.Range("P165").Value = 700
$ws.Range("D166").Value = 44340
$ws.Range("I166").Value = 'Primera'
$ws.Range("J166").Value = 3000
$ws.Range("K166").Value = 600
$ws.Range("L166").Value = 600
$ws.Range("M166").Value = 600
$ws.Range("O166").Value = 'Región del Maule'
$ws.Range("P166").Value = 600
$ws.Range("D167").Value = 44515
$ws.Range("I167").Value = 'Primera'
$ws.Range("J167").Value = 4000
$ws.Range("K167").Value = 500
$ws.Range("L167").Value = 500
$ws.Range("M167").Value = 500
$ws.Range("O167").Value = 'Región del Maule'
$ws.Range("P167").Value = 500
$ws.Range("D168").Value = 44330
$ws.Range("I168").Value = 'Primera'
$ws.Range("J168").Value = 4000
$ws.Range("K168").Value = 650
$ws.Range("L168").Value = 650
$ws.Range("M168").Value = 650
$ws.Range("O168").Value = 'Región del Maule'
$ws.Range("P168").Value = 650
$ws.Range("D169").Value = 44432
$ws.Range("I169").Value = 'Segunda'
$ws.Range("J169").Value = 3000
$ws.Range("K169").Value = 500
$ws.Range("L169").Value = 500
$ws.Range("M169").Value = 500
$ws.Range("O169").Value = 'Región del Maule'
$ws.Range("P169").Value = 500
$ws.Range("D170").Value = 44181
$ws.Range("I170").Value = 'Primera'
$ws.Range("J170").Value = 2000
$ws.Range("K170").Value = 800
$ws.Range("L170").Value = 800
$ws.Range("M170").Value = 800
$ws.Range("O170").Value = 'Región del Maule'
$ws.Range("P170").Value = 800
$ws.Range("D171").Value = 44271
$ws.Range("I171").Value = 'Primera'
$ws.Range("J171").Value = 3000
$ws.Range("K171").Value = 800
$ws.Range("L171").Value = 800
$ws.Range("M171").Value = 800
$ws.Range("O171").Value = 'Región del Maule'
$ws.Range("P171").Value = 800
$ws.Range("D172").Value = 44307
$ws.Range("I172").Value = 'Primera'
$ws.Range("J172").Value = 5000
$ws.Range("K172").Value = 500
$ws.Range("L172").Value = 500
$ws.Range("M172").Value = 500
$ws.Range("O172").Value = 'Región del Maule'
$ws.Range("P172").Value = 500
$ws.Range("D173").Value = 44400
$ws.Range("I173").Value = 'Segunda'
$ws.Range("J173").Value = 3000
$ws.Range("K173").Value = 500
$ws.Range("L173").Value = 500
$ws.Range("M173").Value = 500
$ws.Range("O173").Value = 'Región del Maule'
$ws.Range("P173").Value = 500
$ws.Range("D174").Value = 44309
$ws.Range("I174").Value = 'Primera'
$ws.Range("J174").Value = 5000
$ws.Range("K174").Value = 600
$ws.Range("L174").Value = 600
$ws.Range("M174").Value = 600
$ws.Range("O174").Value = 'Región del Maule'
$ws.Range("P174").Value = 600
$ws.Range("D175").Value = 44508
$ws.Range("I175").Value = 'Primera'
$ws.Range("J175").Value = 5000
$ws.Range("K175").Value = 500
$ws.Range("L175").Value = 500
$ws.Range("M175").Value = 500
$ws.Range("O175").Value = 'Región del Maule'
$ws.Range("P175").Value = 500
